$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.555.76'
$ws.Range('E2').Value = '  +2.79%  '
$ws.Range('D3').Value = '1.669.01'
$ws.Range('E3').Value = '  +2.13%  '
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '237.21'
$ws.Range('E5').Value = '  +1.14%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').Value = '0.4739'
$ws.Range('E7').Value = '  +0.67%  '
$ws.Range('D8').Value = '0.2602'
$ws.Range('E8').Value = '  +1.86%  '
$ws.Range('D9').Value = '0.06171'
$ws.Range('E9').Value = '  +1.82%  '
$ws.Range('D10').Value = '1.668.17'
$ws.Range('E10').Value = '  +2.03%  '
$ws.Range('D11').Value = '0.07004'
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('D12').Value = '14.78'
$ws.Range('E12').Value = '  +1.39%  '
$ws.Range('D13').Value = '0.5864'
$ws.Range('E13').Value = '  -2.98%  '
$ws.Range('D14').Value = '4.364'
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').Value = '75.43'
$ws.Range('E15').Value = '  +3.63%  '
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').Value = '0.9991'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = '25.549.01'
$ws.Range('E18').Value = '  +2.75%  '
$ws.Range('D19').Value = '0.000006728'
$ws.Range('E19').Value = '  +2.55%  '
$ws.Range('D20').Value = '11.41'
$ws.Range('E20').Value = '  +2.74%  '
$ws.Range('D21').Value = '1.882.55'
$ws.Range('E21').Value = '  +1.91%  '
$ws.Range('D22').Value = '4.439'
$ws.Range('E22').Value = '  +2.13%  '
$ws.Range('D23').Value = '8.778'
$ws.Range('E23').Value = '  +2.68%  '
$ws.Range('D24').Value = '5.227'
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').Value = '137.11'
$ws.Range('E25').Value = '  +3.05%  '
$ws.Range('D26').Value = '14.99'
$ws.Range('E26').Value = '  +1.78%  '
$ws.Range('D27').Value = '1.389'
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('D28').Value = '1.719'
$ws.Range('E28').Value = '  +5.72%  '
$ws.Range('D29').Value = '104.39'
$ws.Range('E29').Value = '  +1.24%  '
$ws.Range('D30').Value = '3.999'
$ws.Range('E30').Value = '  +6.08%  '
$ws.Range('D31').Value = '0.07830'
$ws.Range('E31').Value = '  +1.49%  '
$ws.Range('D32').Value = '3.626'
$ws.Range('E32').Value = '  +2.96%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.04309'
$ws.Range('E33').Value = '  +0.51%  '

$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '2.624'
$ws.Range('E34').Value = '  +1.65%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '0.9545'
$ws.Range('E35').Value = '  +3.97%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.6057'
$ws.Range('E36').Value = '  +4.87%  '

$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '0.9368'
$ws.Range('E37').Value = '  +15.35%  '

$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '2.524'
$ws.Range('E38').Value = '  -0.60%  '

$ws.Range('B39').Value = 'PaxDollar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D39').Value = '0.9999'
$ws.Range('E39').Value = '  +0.11%  '

$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '1.854'
$ws.Range('E40').Value = '  +4.87%  '

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.01479'
$ws.Range('E41').Value = '  -3.64%  '

$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Value = '99.36'
$ws.Range('E42').Value = '  +2.39%  '

$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = '0.3749'
$ws.Range('E43').Value = '  +1.94%  '

$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '4.899'
$ws.Range('E44').Value = '  +4.26%  '

$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = '0.1116'
$ws.Range('E45').Value = '  +2.77%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '6.208'
$ws.Range('E46').Value = '  +3.36%  '

$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.05262'
$ws.Range('E47').Value = '  +1.19%  '

$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').Value = '29.87'
$ws.Range('E48').Value = '  +1.69%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '7.484'
$ws.Range('E49').Value = '  +4.79%  '

$ws.Range('B50').Value = 'TrueUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range('D50').Value = '1.002'
$ws.Range('E50').Value = '  +0.15%  '

$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '1.203'
$ws.Range('E51').Value = '  +2.49%  '
